$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "CASOS DE USO" section texts (column E) ---
$ws.Range("E3").Value = "Registrar Preços"
$ws.Range("E4").Value = "Registrar Preços"
$ws.Range("E5").Value = "Manter Cesta de Produtos"
$ws.Range("E10").Value = "Manter Cesta Personalizada"

# --- Update requirement descriptions (column B) ---
$ws.Range("B10").Value = "Permitir que o consumidor crie sua própria cesta de produtos, podendo incluir ou retirar produtos cujos preços deseja acompanhar."
$ws.Range("B10").WrapText = $true

$ws.Range("B12").Value = "Permitir que o consumidor sugira a inclusão ou remoção de produtos na cesta do sistema"

# --- Vertically top-align the "RE Nº" numbers next to wrapped descriptions ---
$ws.Range("A10").VerticalAlignment = -4160
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("A13").VerticalAlignment = -4160

# --- Row 10 now needs extra height to fit the longer wrapped text ---
$ws.Rows("10").RowHeight = 30

# --- Add new "A pensar" notes section below the table ---
$ws.Range("B17").Value = "A pensar:"
$ws.Range("B17").WrapText = $true

$ws.Range("B18").Value = "Imaginar uma forma de validar os preços alimentados pelos pesquisadores, de modo a evitar grandes discrepâncias ou eventuais erros de digitação."
$ws.Range("B18").WrapText = $true
$ws.Rows("18").RowHeight = 45

# --- Update selection to match author's last cursor position ---
$ws.Range("B16").Select()

# --- Page setup for printing ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
